# Apply updated cryptocurrency price/volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.704.32"
$ws.Range("E2").Value = "  -3.45%  "
$ws.Range("D3").Value = "2.094.79"
$ws.Range("E3").Value = "  -2.58%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "344.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5148"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4404"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.77"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09300"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.170"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "2.115.27"
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.288"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.752"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001151"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.009"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("E19").Value = "  +6.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06632"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.007"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.194"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "29.757.41"
$ws.Range("E23").Value = "  -3.56%  "
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.313"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.58%  "
$ws.Range("D26").Value = "2.351.05"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.522"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.130"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1050"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.653"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.172"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.945"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.116"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02568"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06723"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6857"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2228"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.02%  "
$ws.Range("E43").Value = "  +1.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6636"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.316"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.43%  "
$ws.Range("E47").Value = "  -3.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000346"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.221"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3327"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.69%  "
